# Apply updated profit/price figures to each Leve sheet.
# Values below were computed upstream (scheduled data refresh) and
# map 1:1 onto specific cells per sheet (column headers H..N of the
# leve tables): currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = [ordered]@{
    "H76" = 6380.6665
    "I76" = 5723
    "K76" = 5723
    "M76" = -5408
    "H79" = 6380.6665
    "I79" = 5723
    "K79" = 5723
    "M79" = -4631
    "H116" = 5514.1113
    "I116" = 4897.3335
    "J116" = 5822.5
    "K116" = 4897.3335
    "L116" = 5822.5
    "M116" = -1455.3335
    "N116" = -12706.5
    "H137" = 40621.89
    "I137" = 75296.625
    "J137" = 2794.9092
    "K137" = 225889.875
    "L137" = 8384.7276
    "M137" = -223339.875
    "N137" = -13484.7276
    "H138" = 3162.5244
    "I138" = 2720.2666
    "K138" = 8160.7998
    "M138" = -3020.7998
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = [ordered]@{
    "H45" = 6541777
    "J45" = 6731.3335
    "L45" = 6731.3335
    "N45" = -7485.3335
    "H122" = 1392055.4
    "I122" = 3175.6365
    "K122" = 9526.9095
    "M122" = -7076.9095
    "H132" = 2413.842
    "I132" = 1591.6
    "K132" = 4774.799999999999
    "M132" = -2244.799999999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = [ordered]@{
    "H20" = 1535.7742
    "I20" = 1288.5
    "J20" = 1985.3636
    "K20" = 1288.5
    "L20" = 1985.3636
    "M20" = -1041.5
    "N20" = -2479.3636
    "H99" = 14388922
    "I99" = 23978522
    "J99" = 4522.5
    "K99" = 23978522
    "L99" = 4522.5
    "M99" = -23977024
    "N99" = -7518.5
    "H134" = 7186.3076
    "I134" = 3128.375
    "K134" = 9385.125
    "M134" = -6850.125
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = [ordered]@{
    "H31" = 20861.375
    "I31" = 1388.4
    "K31" = 1388.4
    "M31" = -1093.4
    "H34" = 20861.375
    "I34" = 1388.4
    "K34" = 1388.4
    "M34" = -1186.4
    "H53" = 20000
    "J53" = 20000
    "L53" = 20000
    "N53" = -21214
    "H58" = 3107.8809
    "I58" = 3295.5483
    "J58" = 2579
    "K58" = 3295.5483
    "L58" = 2579
    "M58" = -3092.5483
    "N58" = -2985
    "H86" = 13788.206
    "I86" = 12587.046
    "J86" = 15990.333
    "K86" = 12587.046
    "L86" = 15990.333
    "M86" = -11464.046
    "N86" = -18236.333
    "H89" = 13788.206
    "I89" = 12587.046
    "J89" = 15990.333
    "K89" = 62935.23
    "L89" = 79951.66500000001
    "M89" = -57319.23
    "N89" = -91183.66500000001
    "H105" = 1600.0555
    "I105" = 1365.7858
    "K105" = 1365.7858
    "M105" = 381.2141999999999
    "H107" = 1928.5555
    "I107" = 1935.4615
    "J107" = 1910.6
    "K107" = 1935.4615
    "L107" = 1910.6
    "M107" = -15.46149999999989
    "N107" = -5750.6
    "H109" = 34497.168
    "J109" = 38396.6
    "L109" = 38396.6
    "N109" = -40476.6
    "H132" = 49020.527
    "I132" = 1781.6875
    "K132" = 5345.0625
    "M132" = -2815.0625
    "H136" = 3107.8809
    "I136" = 3295.5483
    "J136" = 2579
    "K136" = 9886.644899999999
    "L136" = 7737
    "M136" = -7336.644899999999
    "N136" = -12837
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = [ordered]@{
    "H2" = 269.33334
    "I2" = 19.714285
    "J2" = 453.26315
    "K2" = 118.28571
    "L2" = 2719.5789
    "M2" = -5.285709999999995
    "N2" = -2945.5789
    "H23" = 265.625
    "I23" = 75
    "J23" = 292.85715
    "K23" = 225
    "L23" = 878.5714499999999
    "M23" = 10
    "N23" = -1348.57145
    "H34" = 533
    "I34" = 21.666666
    "J34" = 1300
    "K34" = 64.99999800000001
    "L34" = 3900
    "M34" = 19.00000199999999
    "N34" = -4068
    "H39" = 3724.3333
    "I39" = 1500
    "J39" = 4169.2
    "K39" = 4500
    "L39" = 12507.6
    "M39" = -4206
    "N39" = -13095.6
    "H55" = 61035
    "J55" = 86073.914
    "L55" = 258221.742
    "N55" = -258575.742
    "H56" = 9621042
    "I56" = 9621042
    "K56" = 9621042
    "M56" = -9620512
    "H61" = 174.4
    "I61" = 143
    "K61" = 429
    "M61" = -214
    "H121" = 3331.6667
    "J121" = 4500
    "L121" = 13500
    "N121" = -16120
    "H134" = 3112
    "I134" = 3123.2
    "K134" = 9369.599999999999
    "M134" = -4299.599999999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = [ordered]@{
    "H70" = 4550171.5
    "I70" = 6065005.5
    "J70" = 5669.8184
    "K70" = 6065005.5
    "L70" = 5669.8184
    "M70" = -6064735.5
    "N70" = -6209.8184
    "H73" = 4550171.5
    "I73" = 6065005.5
    "J73" = 5669.8184
    "K73" = 6065005.5
    "L73" = 5669.8184
    "M73" = -6064069.5
    "N73" = -7541.8184
    "H113" = 33335532
    "I113" = 166666670
    "J113" = 2747.25
    "K113" = 166666670
    "L113" = 2747.25
    "M113" = -166664500
    "N113" = -7087.25
    "H122" = 374204.62
    "I122" = 594969.2
    "J122" = 6263.6665
    "K122" = 1784907.6
    "L122" = 18790.9995
    "M122" = -1782457.6
    "N122" = -23690.9995
    "H134" = 36998.5
    "J134" = 36998.5
    "L134" = 110995.5
    "N134" = -116065.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = [ordered]@{
    "H40" = 11036.875
    "I40" = 8575
    "J40" = 13498.75
    "K40" = 8575
    "L40" = 13498.75
    "M40" = -8439
    "N40" = -13770.75
    "H46" = 5642.36
    "J46" = 6074.643
    "L46" = 6074.643
    "N46" = -6450.643
    "H55" = 1769.7273
    "I55" = 1089.1305
    "J55" = 3335.1
    "K55" = 1089.1305
    "L55" = 3335.1
    "M55" = -916.1305
    "N55" = -3681.1
    "H109" = 21995
    "J109" = 21995
    "L109" = 21995
    "N109" = -24769
    "H122" = 5539.0713
    "I122" = 3262.6667
    "K122" = 9788.000100000001
    "M122" = -7338.000100000001
    "H131" = 39671.43
    "J131" = 42783.332
    "L131" = 42783.332
    "N131" = -52863.332
    "H132" = 4221.0586
    "I132" = 3651.6216
    "K132" = 10954.8648
    "M132" = -8424.864799999999
    "H136" = 23412.082
    "I136" = 28275.54
    "J136" = 4444.6
    "K136" = 84826.62
    "L136" = 13333.8
    "M136" = -82276.62
    "N136" = -18433.8
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = [ordered]@{
    "H96" = 2000
    "J96" = 2000
    "L96" = 2000
    "N96" = -4746
    "H109" = 35992.5
    "J109" = 35992.5
    "L109" = 35992.5
    "N109" = -38766.5
    "H126" = 1998.7084
    "J126" = 1613.6364
    "L126" = 4840.9092
    "N126" = -9780.9092
    "H132" = 38498136
    "I132" = 41668780
    "K132" = 125006340
    "M132" = -125003810
    "H137" = 80000
    "I137" = 40000
    "J137" = 100000
    "K137" = 40000
    "L137" = 100000
    "M137" = -34900
    "N137" = -110200
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
